# Apply crypto price/volume updates per commit
# "Updated cryptos list on Sun Jan 14 05:53:46 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price (column D) values are plain decimal numbers (e.g. "303.88") which Excel
# would otherwise auto-detect and store as a numeric cell. The source data keeps these
# as plain text, so force each such cell to Text format before writing the new value,
# then restore its default ("Normal") style afterwards so no stray style index sticks
# around on the cell (only the cell VALUE should change).
$forceTextCells = @(
    "D5",
    "D6",
    "D10",
    "D11",
    "D13",
    "D15",
    "D16",
    "D18",
    "D20",
    "D21",
    "D22",
    "D23",
    "D25",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D35",
    "D36",
    "D37",
    "D38",
    "D40",
    "D41",
    "D45",
    "D46",
    "D49",
    "D50"
)
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "42.802.21"
$ws.Range("E2").Value = "  -0.66%  "

# Row 3
$ws.Range("D3").Value = "2.556.29"
$ws.Range("E3").Value = "  +0.11%  "

# Row 4
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").Value = "303.88"
$ws.Range("E5").Value = "  +1.84%  "

# Row 6
$ws.Range("D6").Value = "98.25"
$ws.Range("E6").Value = "  +6.24%  "

# Row 7
$ws.Range("E7").Value = "  -0.28%  "

# Row 8
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
$ws.Range("E9").Value = "  -0.85%  "

# Row 10
$ws.Range("D10").Value = "36.97"
$ws.Range("E10").Value = "  +2.66%  "

# Row 11
$ws.Range("D11").Value = "0.0809"
$ws.Range("E11").Value = "  -0.12%  "

# Row 12
$ws.Range("E12").Value = "  +7.95%  "

# Row 13
$ws.Range("D13").Value = "7.57"
$ws.Range("E13").Value = "  -2.32%  "

# Row 14
$ws.Range("D14").Value = "2.511.95"
$ws.Range("E14").Value = "  -0.87%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "14.89"
$ws.Range("E15").Value = "  +4.86%  "

# Row 16
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "0.879"
$ws.Range("E16").Value = "  +0.92%  "

# Row 17
$ws.Range("D17").Value = "42.838.51"
$ws.Range("E17").Value = "  -0.54%  "

# Row 18
$ws.Range("D18").Value = "13.27"
$ws.Range("E18").Value = "  +5.73%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0986"
$ws.Range("E19").Value = "  +0.71%  "

# Row 20
$ws.Range("D20").Value = "6.63"
$ws.Range("E20").Value = "  -0.66%  "

# Row 21
$ws.Range("D21").Value = "71.70"
$ws.Range("E21").Value = "  -0.72%  "

# Row 22
$ws.Range("D22").Value = "255.56"
$ws.Range("E22").Value = "  -2.30%  "

# Row 23
$ws.Range("D23").Value = "2.97"
$ws.Range("E23").Value = "  +1.76%  "

# Row 24
$ws.Range("E24").Value = "  -2.22%  "

# Row 25
$ws.Range("D25").Value = "27.87"
$ws.Range("E25").Value = "  -6.37%  "

# Row 26
$ws.Range("E26").Value = "  -0.18%  "

# Row 27
$ws.Range("D27").Value = "10.10"
$ws.Range("E27").Value = "  -0.11%  "

# Row 28
$ws.Range("D28").Value = "38.06"
$ws.Range("E28").Value = "  +3.70%  "

# Row 29
$ws.Range("D29").Value = "2.10"
$ws.Range("E29").Value = "  -1.55%  "

# Row 30
$ws.Range("D30").Value = "6.00"
$ws.Range("E30").Value = "  +0.16%  "

# Row 31
$ws.Range("D31").Value = "156.99"
$ws.Range("E31").Value = "  +1.93%  "

# Row 32
$ws.Range("E32").Value = "  -0.30%  "

# Row 33
$ws.Range("E33").Value = "  +0.59%  "

# Row 34
$ws.Range("E34").Value = "  +1.37%  "

# Row 35
$ws.Range("D35").Value = "3.31"
$ws.Range("E35").Value = "  -2.98%  "

# Row 36
$ws.Range("D36").Value = "26.36"
$ws.Range("E36").Value = "  +8.54%  "

# Row 37
$ws.Range("D37").Value = "18.55"
$ws.Range("E37").Value = "  +10.41%  "

# Row 38
$ws.Range("D38").Value = "0.115"
$ws.Range("E38").Value = "  -1.09%  "

# Row 39
$ws.Range("E39").Value = "  -0.32%  "

# Row 40
$ws.Range("D40").Value = "2.12"
$ws.Range("E40").Value = "  +34.27%  "

# Row 41
$ws.Range("D41").Value = "3.42"
$ws.Range("E41").Value = "  -1.86%  "

# Row 42
$ws.Range("E42").Value = "  +0.41%  "

# Row 43
$ws.Range("E43").Value = "  -2.65%  "

# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.064.98"
$ws.Range("E44").Value = "  -0.91%  "

# Row 45
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  +0.12%  "

# Row 46
$ws.Range("D46").Value = "87.56"
$ws.Range("E46").Value = "  +1.84%  "

# Row 47
$ws.Range("E47").Value = "  +5.59%  "

# Row 48
$ws.Range("D48").Value = "2.803.19"
$ws.Range("E48").Value = "  +0.25%  "

# Row 49
$ws.Range("D49").Value = "74.98"
$ws.Range("E49").Value = "  +7.92%  "

# Row 50
$ws.Range("D50").Value = "103.31"
$ws.Range("E50").Value = "  -1.22%  "

# Row 51
$ws.Range("E51").Value = "  +1.37%  "

# Restore default style on the forced-text cells so only the VALUE changed, not styling
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).Style = "Normal"
}
